$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new key/value pair (data_show_map / true) -----------------
# A3: plain text label, styled like A2 (Consolas family) but a smaller
# 7pt size. Clone A2's format first (keeps a single shared font/style
# instead of forking one per property) and then shrink just the size.
$ws.Range("A3").Value = "data_show_map"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("A3").Font.Size = 7

# B3: must be the literal TEXT "true" (not a boolean). Typing "true"
# directly gets auto-converted to a Boolean by Excel, so build it via a
# text formula first, then convert the formula result down to a plain
# value in place (Copy + PasteSpecial values) so the stored cell is a
# plain shared string, same as every other text cell on the sheet.
$ws.Range("B3").Formula = '="true"'
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# --- Selection: Excel ends up with B3 selected after the edit ---------
$ws.Range("B3").Select() | Out-Null
